# Auto-generated edit script: updates "想去人数" (want-to-go count) values (col F)
# across all four sheets, and marks one event as "暂时售罄" (temporarily sold out) in col G.
$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8266
$ws1.Range("F3").Value = 129
$ws1.Range("F4").Value = 100
$ws1.Range("F5").Value = 35400
$ws1.Range("F6").Value = 49
$ws1.Range("F8").Value = 720
$ws1.Range("F9").Value = 465
$ws1.Range("F10").Value = 149
$ws1.Range("F11").Value = 450
$ws1.Range("F12").Value = 813
$ws1.Range("F13").Value = 67
$ws1.Range("F14").Value = 635
$ws1.Range("F15").Value = 442
$ws1.Range("F17").Value = 579
$ws1.Range("F18").Value = 157
$ws1.Range("F19").Value = 429
$ws1.Range("F20").Value = 426
$ws1.Range("F21").Value = 1121
$ws1.Range("F23").Value = 737
$ws1.Range("F24").Value = 2390
$ws1.Range("F25").Value = 868
$ws1.Range("F26").Value = 506
$ws1.Range("F28").Value = 1102
$ws1.Range("F30").Value = 661
$ws1.Range("F31").Value = 661
$ws1.Range("F32").Value = 14
$ws1.Range("F33").Value = 1097

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 71
$ws2.Range("F4").Value = 355
$ws2.Range("F5").Value = 322
$ws2.Range("F12").Value = 7

# 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 554

# 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 554
$ws4.Range("F3").Value = 8266
$ws4.Range("F4").Value = 129
$ws4.Range("F5").Value = 100
$ws4.Range("F7").Value = 35400
$ws4.Range("F8").Value = 49
$ws4.Range("F10").Value = 720
$ws4.Range("F11").Value = 465
$ws4.Range("F12").Value = 71
$ws4.Range("F13").Value = 149
$ws4.Range("F14").Value = 450
$ws4.Range("F15").Value = 355
$ws4.Range("F16").Value = 322
$ws4.Range("F18").Value = 813
$ws4.Range("F19").Value = 67
$ws4.Range("F20").Value = 635
$ws4.Range("F21").Value = 442
$ws4.Range("F28").Value = 579
$ws4.Range("F29").Value = 157
$ws4.Range("F30").Value = 429
$ws4.Range("F31").Value = 426
$ws4.Range("F32").Value = 1121
$ws4.Range("F34").Value = 737
$ws4.Range("F35").Value = 2390
$ws4.Range("F36").Value = 868
$ws4.Range("F37").Value = 506
$ws4.Range("F39").Value = 1102
$ws4.Range("F41").Value = 7
$ws4.Range("F42").Value = 661
$ws4.Range("F43").Value = 661
$ws4.Range("F44").Value = 14
$ws4.Range("F45").Value = 1097

# COMICUP 2024SP in 全部类型 sheet is now temporarily sold out: price cell becomes text
$ws4.Range("G7").Value = "暂时售罄"

